# "Changes in loginPage Event Class To check err"
# Sheet1 (login page end-user data) content is cleared so the login page
# test can exercise its error-checking path, and Sheet1 becomes the
# active/selected tab (previously Sheet2 was active).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Clear out the header/data values on Sheet1 (A1:D2) but keep the cell
# formatting (styles) intact.
$ws1.Range("A1:D2").ClearContents()

# Remove the now-orphaned hyperlink that pointed at the (now blank) D2 cell.
$ws1.Hyperlinks.Delete()

# Sheet1 becomes the selected/active sheet (it previously was Sheet2).
$ws1.Activate()
